$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Auto-detect starting row: insert 2 blank rows above the existing table on Sheet1,
# shifting all data down so the table now starts at row 3.
$ws1.Range("A1:D2").Insert()

# Restore Sheet1's remembered selection after the insert.
$ws1.Range("G16").Select() | Out-Null

# Append a new data row to Sheet2.
$ws2.Range("A12").Value = 31
$ws2.Range("B12").Value = "May"
$ws2.Range("C12").Value = "Bach"
$ws2.Range("D12").Value = "Laos"

# Sheet2 becomes the active sheet, with its own remembered selection.
$ws2.Activate()
$ws2.Range("D13").Select() | Out-Null
